$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "locked_out_user" row (row 3) and shift the rows below it up.
$ws.Rows.Item(3).Delete()

# Update the selection to match the target state.
$ws.Range("C11").Select()
